# Signup Mortgage and Login Api's - update Endpoint Definition sheet
# (Request/Response contracts for "Register Mortgage" and "Login" rows,
#  and clear the now-unused contracts on "Transaction Details" / "Admin Approval")

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Endpoint Definition")
$ws2.Activate()

# Row 3 - Register Mortgage: Response Contract now includes login/mortgage fields
$ws2.Range("F3").Value2 = "{loginId:String,`n password:String,`n mortgageNumber:String,`naccountNumber:String,`nmessage:String`n}"

# Row 4 - Login: Request Contract becomes an empty object, Response Contract becomes a simple message
$ws2.Range("E4").Value2 = "{`n}"
$ws2.Range("F4").Value2 = "{message:String`n}`n"
$ws2.Rows.Item(4).RowHeight = 38.25

# Row 6 - Transaction Details: Response Contract cleared (no longer used here)
$ws2.Range("F6").Value2 = ""
$ws2.Rows.Item(6).AutoFit()

# Row 7 - Admin Approval: Response Contract cleared (no longer used here)
$ws2.Range("F7").Value2 = ""
$ws2.Rows.Item(7).AutoFit()

# Update the view: scroll down a bit and select F7 as the active cell
$ws2.Range("F7").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
